$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.228
$ws.Range("E4").Value = 13.168
$ws.Range("E6").Value = 13.024
$ws.Range("A9").Value = -20.912
$ws.Range("E10").Value = 12.519
$ws.Range("B11").Value = 6.528999999999999
$ws.Range("E11").Value = 12.254
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.828
$ws.Range("C21").Value = -12.282
$ws.Range("E21").Value = 13.188
